$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First paragraph of the document ("Loading Mode" heading): right
#    align it by adding <w:jc w:val="right"/> to its paragraph
#    properties.
# ---------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$firstPara.Format.Alignment = 2   # wdAlignParagraphRight

# ---------------------------------------------------------------------
# 2) The "option = { Fill:Element,status:boolean};" paragraph (the one
#    with the 1620/540 indent, i.e. the *second* "Fill:Element"
#    occurrence in the document): right align it, then split the
#    "Fill:Element" run into "Fill:" + "Element" and move the
#    "_GoBack" bookmark so that it wraps the new "Element" run.
# ---------------------------------------------------------------------

# Locate the first "Fill:Element" occurrence so we can search past it
# for the second one (the target paragraph).
$rngAll = $d.Content
$firstHit = $rngAll.Find.Execute("Fill:Element", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Search the remainder of the document for the second occurrence.
$rngRest = $d.Range($rngAll.End, $d.Content.End)
$found = $rngRest.Find.Execute("Fill:Element", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Right-align the paragraph that contains this match.
$targetPara = $rngRest.Paragraphs(1)
$targetPara.Format.Alignment = 2   # wdAlignParagraphRight

# $rngRest now spans exactly "Fill:Element". "Fill:" is the first 5
# characters, "Element" is the remaining 7, so the bookmark should
# wrap just the "Element" part.
$elementRange = $d.Range($rngRest.Start + 5, $rngRest.End)

# Re-adding the "_GoBack" bookmark at the new location both creates
# the bookmarkStart/bookmarkEnd pair around "Element" (splitting the
# run into "Fill:" and "Element") and removes it from its old location
# further down in the document (a bookmark name is unique within the
# doc).
$d.Bookmarks.Add("_GoBack", $elementRange)
